$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.902.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.34%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.503.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.10%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5: BNB -> BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.60"

# Row 6: Solana -> Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.14"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.13%  "

# Row 7: USDC -> USDC
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8: XRP -> XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.68%  "

# Row 9: Dogecoin -> Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.142"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.30%  "

# Row 10: TRON -> TRON
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.08%  "

# Row 11: Cardano -> Cardano
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.58%  "

# Row 12: Toncoin -> Toncoin
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.18%  "

# Row 13: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.962.05"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.48%  "

# Row 14: Avalanche -> Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.77"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.78%  "

# Row 15: WrappedBTC -> WrappedBTC
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.804.21"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.25%  "

# Row 16: ShibaInu -> ShibaInu
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.13%  "

# Row 17: WrappedEther -> WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.501.28"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.89%  "

# Row 18: Chainlink -> Uniswap
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.54"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.15%  "

# Row 19: Uniswap -> Chainlink
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.97"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.39%  "

# Row 20: BitcoinCash -> BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "352.26"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.90%  "

# Row 21: Polkadot -> Polkadot
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.88%  "

# Row 22: Dai -> Litecoin
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.32"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.67%  "

# Row 23: Litecoin -> Dai
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.01%  "

# Row 24: NEARProtocol -> NEARProtocol
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.27"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.47%  "

# Row 25: SuiNetwork -> SuiNetwork
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.74"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.28%  "

# Row 26: Aptos -> Aptos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.14"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.38%  "

# Row 27: WrappedeETH -> WrappedeETH
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.632.88"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.00%  "

# Row 28: Binance-PegBSC-USD -> Binance-PegBSC-USD
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.30%  "

# Row 29: PEPE -> PEPE
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0912"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.29%  "

# Row 30: Bittensor -> Bittensor
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "510.94"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.93%  "

# Row 31: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.29%  "

# Row 32: Fetch.AI -> Fetch.AI
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.28%  "

# Row 33: PancakeSwap -> PancakeSwap
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.92%  "

# Row 34: FirstDigitalUSD -> FirstDigitalUSD
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.02%  "

# Row 35: Kaspa -> Kaspa
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.94%  "

# Row 36: Monero -> Monero
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.20"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.76%  "

# Row 37: EthereumClassic -> EthereumClassic
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.39"

# Row 38: WhiteBITCoin -> WhiteBITCoin
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.63"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.35%  "

# Row 39: ImmutableX -> ImmutableX
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.33"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.02%  "

# Row 40: USDe -> USDe
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.04%  "

# Row 41: Stacks -> Stacks
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.73"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.91%  "

# Row 42: PolygonEcosystemToken -> PolygonEcosystemToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.330"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.66%  "

# Row 43: RenderToken -> RenderToken
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.00%  "

# Row 44: dogwifhat -> dogwifhat
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.40%  "

# Row 45: Aave -> Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.28"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.75%  "

# Row 46: Filecoin -> Filecoin
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.69%  "

# Row 47: BabyDogeCoin -> ARBITRUM
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.519"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.45%  "

# Row 48: ARBITRUM -> Cronos
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0741"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.84%  "

# Row 49: Cronos -> Optimism
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.59"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.55%  "

# Row 50: Optimism -> Mantle
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.583"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.31%  "

# Row 51: Mantle -> Fantom
$ws.Range("B51").Value = "Fantom"
$ws.Range("C51").Value = "https://coinranking.com/coin/uIEWfMFnQo9K_+fantom-ftm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.683"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.06%  "
